$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A14").Value = "24.9.2025"
$ws.Range("B14").Value = 0.375
$ws.Range("C14").Value = 0.52083333333333337

$ws.Range("A15").Value = "25.9.2025"
$ws.Range("B15").Value = 0.64583333333333337
$ws.Range("C15").Value = 0.78125

$ws.Range("B14:C15").NumberFormat = "h:mm AM/PM"

$ws.Range("D18").Select()
